$d = $word.ActiveDocument

$replacements = @(
    @{old="560×7="; new="919×8="},
    @{old="402×9="; new="506×9="},
    @{old="458×9="; new="528×8="},
    @{old="378×5="; new="737×9="},
    @{old="142×5="; new="254×2="},
    @{old="197×9="; new="453×8="},
    @{old="822×8="; new="423×2="},
    @{old="890×9="; new="133×3="},
    @{old="984×3="; new="207×3="},
    @{old="726×5="; new="123×5="},
    @{old="414×5="; new="157×5="},
    @{old="701×2="; new="323×5="},
    @{old="742×4="; new="533×3="},
    @{old="954×5="; new="826×3="},
    @{old="143×9="; new="481×8="},
    @{old="713×8="; new="303×4="},
    @{old="977×8="; new="325×6="},
    @{old="767×8="; new="508×6="},
    @{old="146×4="; new="463×6="},
    @{old="293×9="; new="165×4="},
    @{old="954×3="; new="207×4="},
    @{old="533×8="; new="968×4="},
    @{old="922×5="; new="500×7="},
    @{old="973×3="; new="357×9="},
    @{old="555×4="; new="868×5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
